# The author's commit ("wrapping up test file audit") re-saved this
# workbook from a newer Excel build while finishing an audit pass. The
# meaningful, content-level edit buried in that resave is: the stray
# leftover row 16 ("Sheet" / 3 / 4) on the "optimization_parameters"
# sheet was selected and deleted (Excel then shifts the old row 17 up
# into row 16). Focus also ends up on "optimization_diagnostics" (the
# last sheet), matching the workbook's new active-tab/selection state.

$wb = $excel.ActiveWorkbook

# Go to the sheet with the orphaned row, select the whole row the way a
# user would before deleting it, then delete it.
$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
[void]$paramsSheet.Activate()
[void]$paramsSheet.Rows(16).Select()
[void]$paramsSheet.Rows(16).Delete()

# Finish up on the diagnostics sheet (now the active/selected tab).
$diagSheet = $wb.Worksheets.Item("optimization_diagnostics")
[void]$diagSheet.Activate()
